$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("docs")
$ws.Range("B68").Value = "What I Believe - Einstein"
$ws.Range("A68").Value = "what-i-believe"
$ws.Range("C68").Value = "philosophy"
$ws.Range("D68").Value = "None"
$ws.Range("D68").Font.Color = 0
$ws.Range("D68").Select()
